$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header label in row 2, column B (was "unnamed: 1_level_1" -> "total")
$ws.Range("B2").Value = "total"

# Shift data rows up: row6 values -> row5, row7 -> row6, row9 -> row7,
# row10 -> row8, row11 -> row9, row12 -> row10, row13 -> row11.
# The row labels in column A need to move along with the data, since
# the "situação do domicílio" / "grandes regiões" sub-header rows are
# being removed.

$ws.Range("A5").Value = "urbana"
$ws.Range("A6").Value = "rural"
$ws.Range("A7").Value = "norte"
$ws.Range("A8").Value = "nordeste"
$ws.Range("A9").Value = "sudeste"
$ws.Range("A10").Value = "sul"
$ws.Range("A11").Value = "centro-oeste"

# Row 5 (urbana) gets the values that used to sit on row 6
$ws.Range("B5").Value = 2.23
$ws.Range("C5").Value = 5.56
$ws.Range("D5").Value = 4.22
$ws.Range("E5").Value = 3.79
$ws.Range("F5").Value = 3.73
$ws.Range("G5").Value = 5.2

# Row 6 (rural) gets the values that used to sit on row 7
$ws.Range("B6").Value = 7.02
$ws.Range("C6").Value = 11.26
$ws.Range("D6").Value = 11.72
$ws.Range("E6").Value = 11.41
$ws.Range("F6").Value = 14.86
$ws.Range("G6").Value = 27.62

# Row 7 (norte) gets the values that used to sit on row 9
$ws.Range("B7").Value = 5.08
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 8.83
$ws.Range("E7").Value = 8.5
$ws.Range("F7").Value = 11.13
$ws.Range("G7").Value = 12.18

# Row 8 (nordeste) gets the values that used to sit on row 10
$ws.Range("B8").Value = 4.09
$ws.Range("C8").Value = 7.79
$ws.Range("D8").Value = 6.31
$ws.Range("E8").Value = 6.16
$ws.Range("F8").Value = 7.3
$ws.Range("G8").Value = 9.69

# Row 9 (sudeste) gets the values that used to sit on row 11
$ws.Range("B9").Value = 3.96
$ws.Range("C9").Value = 11.44
$ws.Range("D9").Value = 9.6
$ws.Range("E9").Value = 7.18
$ws.Range("F9").Value = 5.95
$ws.Range("G9").Value = 8.27

# Row 10 (sul) gets the values that used to sit on row 12
$ws.Range("B10").Value = 4.8
$ws.Range("C10").Value = 14.27
$ws.Range("D10").Value = 9.880000000000001
$ws.Range("E10").Value = 8.31
$ws.Range("F10").Value = 8.17
$ws.Range("G10").Value = 10.45

# Row 11 (centro-oeste) gets the values that used to sit on row 13
$ws.Range("B11").Value = 6.19
$ws.Range("C11").Value = 16.6
$ws.Range("D11").Value = 14
$ws.Range("E11").Value = 7.92
$ws.Range("F11").Value = 9.81
$ws.Range("G11").Value = 12.99

# Remove the now-obsolete trailing rows (old row12 "sul", row13
# "centro-oeste" header duplicates, and row14 source footnote).
$ws.Rows("12:14").Delete()
